# LogTime.xlsx edit script
# - Adds K6=3 on ANLT (sheet1)
# - Adds two new rows (HOME / DETAILS) on ANHDT (sheet3)
# - Updates the active sheet / selections to match the authored state
#   (QUANGD -> ANHDT becomes the active/selected tab)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # ANLT
$ws2 = $wb.Worksheets.Item(2)   # QUANGD
$ws3 = $wb.Worksheets.Item(3)   # ANHDT

# --- ANLT (sheet1): new data point K6 = 3 ---
$ws1.Range("K6").Value = 3

# --- ANHDT (sheet3): two new rows of data ---
$ws3.Range("A4").Value = "HOME "
$ws3.Range("N4").Value = 4
$ws3.Range("A5").Value = "DETAILS "
$ws3.Range("N5").Value = 4

# --- View / selection state ---
# ANLT: selection moves to J14, no longer the active tab
$ws1.Activate()
$ws1.Range("J14").Select()

# QUANGD: selection becomes B2:L10, no longer the active tab
$ws2.Activate()
$ws2.Range("B2:L10").Select()

# ANHDT: becomes the active tab, scrolled so column B is leftmost,
# with P10 selected
$ws3.Activate()
$ws3.Range("P10").Select()
